$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 249.66667
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 249.66667
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 249.66667
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -475.66667
$ws.Range("H17").Value = 1716.5555
$ws.Range("J17").Value = 1716.5555
$ws.Range("L17").Value = 5149.666499999999
$ws.Range("N17").Value = -5485.666499999999
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()
$ws.Range("H51").Value = 10351.667
$ws.Range("J51").Value = 10555.5
$ws.Range("L51").Value = 10555.5
$ws.Range("N51").Value = -11523.5
$ws.Range("H53").Value = 704.6667
$ws.Range("I53").Value = 542.75
$ws.Range("K53").Value = 542.75
$ws.Range("M53").Value = 94.25
$ws.Range("H107").Value = 537.6667
$ws.Range("J107").Value = 1005.3333
$ws.Range("L107").Value = 1005.3333
$ws.Range("N107").Value = -4845.3333
$ws.Range("H116").Value = 14697
$ws.Range("I116").Value = 3995
$ws.Range("K116").Value = 3995
$ws.Range("M116").Value = -553

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 272.5
$ws.Range("I2").Value = 272.5
$ws.Range("K2").Value = 272.5
$ws.Range("M2").Value = -159.5
$ws.Range("H92").Value = 70000
$ws.Range("J92").Value = 50000
$ws.Range("L92").Value = 50000
$ws.Range("N92").Value = -54992
$ws.Range("H116").Value = 272.5
$ws.Range("I116").Value = 272.5
$ws.Range("K116").Value = 272.5
$ws.Range("M116").Value = 2021.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 272.5
$ws.Range("I3").Value = 272.5
$ws.Range("K3").Value = 272.5
$ws.Range("M3").Value = -158.5
$ws.Range("H5").Value = 847.0769
$ws.Range("I5").Value = 238.33333
$ws.Range("K5").Value = 238.33333
$ws.Range("M5").Value = -125.33333
$ws.Range("H7").Value = 25050500
$ws.Range("I7").Value = 25050500
$ws.Range("K7").Value = 25050500
$ws.Range("M7").Value = -25050387
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 200
$ws.Range("K22").Value = 200
$ws.Range("M22").Value = -27

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2318
$ws.Range("I16").Value = 4895
$ws.Range("J16").Value = 600
$ws.Range("K16").Value = 4895
$ws.Range("L16").Value = 600
$ws.Range("M16").Value = -4608
$ws.Range("N16").Value = -1174
$ws.Range("H113").Value = 2318
$ws.Range("I113").Value = 4895
$ws.Range("J113").Value = 600
$ws.Range("K113").Value = 4895
$ws.Range("L113").Value = 600
$ws.Range("M113").Value = -2725
$ws.Range("N113").Value = -4940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 144.18182
$ws.Range("I17").Value = 70.5
$ws.Range("J17").Value = 340.66666
$ws.Range("K17").Value = 211.5
$ws.Range("L17").Value = 1021.99998
$ws.Range("M17").Value = -42.5
$ws.Range("N17").Value = -1359.99998
$ws.Range("H23").Value = 213.23077
$ws.Range("I23").Value = 157.6
$ws.Range("J23").Value = 248
$ws.Range("K23").Value = 472.8
$ws.Range("L23").Value = 744
$ws.Range("M23").Value = -237.8
$ws.Range("N23").Value = -1214
$ws.Range("H33").Value = 355
$ws.Range("I33").Value = 200
$ws.Range("J33").Value = 510
$ws.Range("K33").Value = 1200
$ws.Range("L33").Value = 3060
$ws.Range("M33").Value = -917
$ws.Range("N33").Value = -3626
$ws.Range("H40").Value = 47.77778
$ws.Range("J40").Value = 51.25
$ws.Range("L40").Value = 205
$ws.Range("N40").Value = -343
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H80").Value = 3333
$ws.Range("I80").Value = 3333
$ws.Range("K80").Value = 9999
$ws.Range("M80").Value = -9063
$ws.Range("H83").Value = 3333
$ws.Range("I83").Value = 3333
$ws.Range("K83").Value = 29997
$ws.Range("M83").Value = -25317
$ws.Range("H88").Value = 15000
$ws.Range("J88").Value = 15000
$ws.Range("L88").Value = 45000
$ws.Range("N88").Value = -45856
$ws.Range("H91").Value = 15000
$ws.Range("J91").Value = 15000
$ws.Range("L91").Value = 45000
$ws.Range("N91").Value = -47964
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("M108").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 5004
$ws.Range("I4").Value = 5004
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 5004
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -4892
$ws.Range("N4").ClearContents()
$ws.Range("H70").Value = 6626.857
$ws.Range("I70").Value = 6908
$ws.Range("K70").Value = 6908
$ws.Range("M70").Value = -6638
$ws.Range("H73").Value = 6626.857
$ws.Range("I73").Value = 6908
$ws.Range("K73").Value = 6908
$ws.Range("M73").Value = -5972
$ws.Range("H122").Value = 1389.8
$ws.Range("I122").Value = 1389.8
$ws.Range("K122").Value = 4169.4
$ws.Range("M122").Value = -1719.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 25500
$ws.Range("J12").Value = 50000
$ws.Range("L12").Value = 50000
$ws.Range("N12").Value = -50340
$ws.Range("H40").Value = 5000
$ws.Range("J40").Value = 5000
$ws.Range("L40").Value = 5000
$ws.Range("N40").Value = -5272
$ws.Range("H46").Value = 1000
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H55").Value = 421.42856
$ws.Range("I55").Value = 375
$ws.Range("J55").Value = 483.33334
$ws.Range("K55").Value = 375
$ws.Range("L55").Value = 483.33334
$ws.Range("M55").Value = -202
$ws.Range("N55").Value = -829.33334
$ws.Range("H58").Value = 50000
$ws.Range("J58").Value = 50000
$ws.Range("L58").Value = 50000
$ws.Range("N58").Value = -50520

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4992
$ws.Range("I62").Value = 4992
$ws.Range("K62").Value = 4992
$ws.Range("M62").Value = -4368
$ws.Range("H65").Value = 4992
$ws.Range("I65").Value = 4992
$ws.Range("K65").Value = 24960
$ws.Range("M65").Value = -21840
$ws.Range("H126").Value = 4748.5
$ws.Range("I126").Value = 4499.5
$ws.Range("K126").Value = 13498.5
$ws.Range("M126").Value = -11028.5
